# Generate Report for Handoff
# Updates the "Latest Handoff"/"Latest HO Xliff Generate" timestamps produced by the
# handoff report generation run, and sets the Priority column to "ht" for the rows
# that were (re)generated in this handoff pass.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$rows = @(7, 8, 9, 10, 11, 14)

foreach ($r in $rows) {
    # Overview sheet: "Latest HO Xliff Generate Date" column (G)
    $overview.Range("G$r").Value = "2016-08-20 04:19:29"

    # zh-cn sheet: "Latest Handoff Datetime" column (H) and "Priority" column (E)
    $zhcn.Range("H$r").Value = "2016-08-20 04:19:24"
    $zhcn.Range("E$r").Value = "ht"

    # de-de sheet: "Latest Handoff Datetime" column (H) and "Priority" column (E)
    $dede.Range("H$r").Value = "2016-08-20 04:19:29"
    $dede.Range("E$r").Value = "ht"
}
